$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(136, 1).Value = 135
$ws.Range("B136").Value = 'Friday, Jan 13'
$ws.Range("C136").Value = '2:50 PM'
$ws.Range("D136").Value = 'FR1937'
$ws.Range("E136").Value = 'Helsinki'
$ws.Range("F136").Value = '(HEL)'
$ws.Range("G136").Value = 'Ryanair '
$ws.Range("H136").Value = 'B738'
$ws.Range("I136").Value = '(SP-RKW)'
$ws.Range("J136").Value = '2:38 PM'
$ws.Range("L136").Value = '0 hours, -12 minutes'

$ws.Cells.Item(137, 1).Value = 136
$ws.Range("B137").Value = 'Friday, Jan 13'
$ws.Range("C137").Value = '3:15 PM'
$ws.Range("D137").Value = 'FR2669'
$ws.Range("E137").Value = 'London'
$ws.Range("F137").Value = '(STN)'
$ws.Range("G137").Value = 'Ryanair '
$ws.Range("H137").Value = 'B738'
$ws.Range("I137").Value = '(EI-DWC)'
$ws.Range("J137").Value = '3:01 PM'
$ws.Range("L137").Value = '0 hours, -14 minutes'

$ws.Cells.Item(138, 1).Value = 137
$ws.Range("B138").Value = 'Friday, Jan 13'
$ws.Range("C138").Value = '3:25 PM'
$ws.Range("D138").Value = 'FR2264'
$ws.Range("E138").Value = 'Lisbon'
$ws.Range("F138").Value = '(LIS)'
$ws.Range("G138").Value = 'Ryanair '
$ws.Range("H138").Value = 'B738'
$ws.Range("I138").Value = '(SP-RKP)'
$ws.Range("J138").Value = '4:02 PM'
$ws.Range("L138").Value = '0 hours, 37 minutes'

$ws.Cells.Item(139, 1).Value = 138
$ws.Range("B139").Value = 'Friday, Jan 13'
$ws.Range("C139").Value = '4:20 PM'
$ws.Range("D139").Value = 'FR9258'
$ws.Range("E139").Value = 'Malta'
$ws.Range("F139").Value = '(MLA)'
$ws.Range("G139").Value = 'Ryanair '
$ws.Range("H139").Value = 'B738'
$ws.Range("I139").Value = '(9H-QAD)'
$ws.Range("J139").Value = '4:12 PM'
$ws.Range("L139").Value = '0 hours, -8 minutes'

$ws.Cells.Item(140, 1).Value = 139
$ws.Range("B140").Value = 'Friday, Jan 13'
$ws.Range("C140").Value = '4:45 PM'
$ws.Range("D140").Value = 'FR1112'
$ws.Range("E140").Value = 'Rome'
$ws.Range("F140").Value = '(CIA)'
$ws.Range("G140").Value = 'Ryanair '
$ws.Range("H140").Value = 'B38M'
$ws.Range("I140").Value = '(9H-VUH)'
$ws.Range("J140").Value = '4:15 PM'
$ws.Range("L140").Value = '0 hours, -30 minutes'

$ws.Cells.Item(141, 1).Value = 140
$ws.Range("B141").Value = 'Friday, Jan 13'
$ws.Range("C141").Value = '5:20 PM'
$ws.Range("D141").Value = 'FR1945'
$ws.Range("E141").Value = 'Stockholm'
$ws.Range("F141").Value = '(ARN)'
$ws.Range("G141").Value = 'Buzz '
$ws.Range("H141").Value = 'B38M'
$ws.Range("I141").Value = '(SP-RZE)'
$ws.Range("J141").Value = '5:08 PM'
$ws.Range("L141").Value = '0 hours, -12 minutes'

$ws.Cells.Item(142, 1).Value = 141
$ws.Range("B142").Value = 'Friday, Jan 13'
$ws.Range("C142").Value = '5:50 PM'
$ws.Range("D142").Value = 'FR1943'
$ws.Range("E142").Value = 'Bologna'
$ws.Range("F142").Value = '(BLQ)'
$ws.Range("G142").Value = 'Ryanair '
$ws.Range("H142").Value = 'B738'
$ws.Range("I142").Value = '(SP-RKD)'
$ws.Range("J142").Value = '5:39 PM'
$ws.Range("L142").Value = '0 hours, -11 minutes'

$ws.Cells.Item(143, 1).Value = 142
$ws.Range("B143").Value = 'Friday, Jan 13'
$ws.Range("C143").Value = '8:55 PM'
$ws.Range("D143").Value = 'FR4534'
$ws.Range("E143").Value = 'Porto'
$ws.Range("F143").Value = '(OPO)'
$ws.Range("G143").Value = 'Buzz '
$ws.Range("H143").Value = 'B38M'
$ws.Range("I143").Value = '(SP-RZG)'
$ws.Range("J143").Value = '8:56 PM'
$ws.Range("L143").Value = '0 hours, 1 minutes'

$ws.Cells.Item(144, 1).Value = 143
$ws.Range("B144").Value = 'Friday, Jan 13'
$ws.Range("C144").Value = '9:10 PM'
$ws.Range("D144").Value = 'FR4554'
$ws.Range("E144").Value = 'Naples'
$ws.Range("F144").Value = '(NAP)'
$ws.Range("G144").Value = 'Ryanair '
$ws.Range("H144").Value = 'B738'
$ws.Range("I144").Value = '(SP-RKW)'
$ws.Range("J144").Value = '9:09 PM'
$ws.Range("L144").Value = '0 hours, -1 minutes'

$ws.Cells.Item(145, 1).Value = 144
$ws.Range("B145").Value = 'Friday, Jan 13'
$ws.Range("C145").Value = '9:15 PM'
$ws.Range("D145").Value = 'FR1021'
$ws.Range("E145").Value = 'London'
$ws.Range("F145").Value = '(STN)'
$ws.Range("G145").Value = 'Ryanair '
$ws.Range("H145").Value = 'B738'
$ws.Range("I145").Value = '(EI-EVC)'
$ws.Range("J145").Value = '9:05 PM'
$ws.Range("L145").Value = '0 hours, -10 minutes'

$ws.Cells.Item(146, 1).Value = 145
$ws.Range("B146").Value = 'Friday, Jan 13'
$ws.Range("C146").Value = '9:15 PM'
$ws.Range("D146").Value = 'FR1107'
$ws.Range("E146").Value = 'Alicante'
$ws.Range("F146").Value = '(ALC)'
$ws.Range("G146").Value = 'Ryanair '
$ws.Range("H146").Value = 'B738'
$ws.Range("I146").Value = '(SP-RSS)'
$ws.Range("J146").Value = '9:33 PM'
$ws.Range("L146").Value = '0 hours, 18 minutes'

$ws.Cells.Item(147, 1).Value = 146
$ws.Range("B147").Value = 'Friday, Jan 13'
$ws.Range("C147").Value = '9:30 PM'
$ws.Range("D147").Value = 'FR2007'
$ws.Range("E147").Value = 'Cologne'
$ws.Range("F147").Value = '(CGN)'
$ws.Range("G147").Value = 'Ryanair '
$ws.Range("H147").Value = 'B738'
$ws.Range("I147").Value = '(9H-QBA)'
$ws.Range("J147").Value = '9:12 PM'
$ws.Range("L147").Value = '0 hours, -18 minutes'

$ws.Cells.Item(148, 1).Value = 147
$ws.Range("B148").Value = 'Friday, Jan 13'
$ws.Range("C148").Value = '9:35 PM'
$ws.Range("D148").Value = 'FR3898'
$ws.Range("E148").Value = 'Milan'
$ws.Range("F148").Value = '(BGY)'
$ws.Range("G148").Value = 'Ryanair '
$ws.Range("H148").Value = 'B738'
$ws.Range("I148").Value = '(SP-RKP)'
$ws.Range("J148").Value = '9:17 PM'
$ws.Range("L148").Value = '0 hours, -18 minutes'

$ws.Cells.Item(149, 1).Value = 148
$ws.Range("B149").Value = 'Friday, Jan 13'
$ws.Range("C149").Value = '11:00 PM'
$ws.Range("D149").Value = 'FR4238'
$ws.Range("E149").Value = 'Bari'
$ws.Range("F149").Value = '(BRI)'
$ws.Range("G149").Value = 'Ryanair '
$ws.Range("H149").Value = 'B738'
$ws.Range("I149").Value = '(SP-RKD)'
$ws.Range("J149").Value = '10:45 PM'
$ws.Range("L149").Value = '0 hours, -15 minutes'
